# Apply the "update questions to add GPS type" edit to the survey sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# 1) Row 11 ("store_gps") changes its question type from plain "text" to "gps_coord".
$ws.Range("A11").Value = "gps_coord"

# 2) Insert 6 new rows before the current row 172 (the "new_section" row that
#    separates the "who do you ... with" questions from the closing
#    questions/notes rows). This pushes everything from the old row 172
#    onward down by 6, matching the diff's row 175-180 tail.
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()

# 3) Rewrite rows 169-174 as the new "select_one yes_no_2" symptom questions.
$ws.Range("A169").Value = "select_one yes_no_2"
$ws.Range("B169").Value = "has_facial_swelling"
$ws.Range("C169").Value = "Do you have facial swelling?"

$ws.Range("A170").Value = "select_one yes_no_2"
$ws.Range("B170").Value = "has_muscle_fatigue"
$ws.Range("C170").Value = "Do you have muscle fatigue?"

$ws.Range("A171").Value = "select_one yes_no_2"
$ws.Range("B171").Value = "has_vomiting"
$ws.Range("C171").Value = "Are you vomiting?"

$ws.Range("A172").Value = "select_one yes_no_2"
$ws.Range("B172").Value = "has_cough"
$ws.Range("C172").Value = "Do you have a cough?"

$ws.Range("A173").Value = "select_one yes_no_2"
$ws.Range("B173").Value = "has_meningitis"
$ws.Range("C173").Value = "Do you have meningitis?"

$ws.Range("A174").Value = "select_one yes_no_2"
$ws.Range("B174").Value = "has_hypertension"
$ws.Range("C174").Value = "Do you have hypertension?"

# 4) The previously-existing "who ... with" / new_section / questions / notes
#    rows now live 6 rows further down (175-180). Restore their original
#    content there (row-insert duplicated formatting/style but left the
#    shifted cells otherwise as before, so just fill the text back in).
$ws.Range("A175").Value = "text"
$ws.Range("B175").Value = "who_live_with"
$ws.Range("C175").Value = "Who do you live with?"

$ws.Range("A176").Value = "text"
$ws.Range("B176").Value = "who_sharefood_with"
$ws.Range("C176").Value = "Who do you share food with?"

$ws.Range("A177").Value = "text"
$ws.Range("B177").Value = "who_work_with"
$ws.Range("C177").Value = "Who do you work with?"

$ws.Range("A178").Value = "new_section"

$ws.Range("A179").Value = "text"
$ws.Range("B179").Value = "questions"
$ws.Range("C179").Value = "Please, do you have any question for me?"

$ws.Range("A180").Value = "text"
$ws.Range("B180").Value = "notes"
$ws.Range("C180").Value = "Include any notes about this interview"

# 5) Update the view: the user scrolled/selected further down the now-longer
#    sheet (selection on C174).
$ws.Range("C174").Select()
